# Add a new bulleted/numbered list item (same list as the existing two
# hyperlink entries, numId=2) at the end of the document, containing the
# new HackerNoon Flask/iOS article URL as plain text.

$d = $word.ActiveDocument

# Start a brand new paragraph right after the current last paragraph
# (the GameKit matchmaking list item).
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

# The newly created (currently empty) paragraph is now the last one;
# give it the exact paragraph/run formatting used by the other list
# items (Body style, list numbering numId=2 at level 0, LTR) and fill
# in the new URL text.
$newPara = $d.Paragraphs.Last
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr>' +
    '<w:pStyle w:val="Body"/>' +
    '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr>' +
    '<w:bidi w:val="0"/>' +
  '</w:pPr>' +
  '<w:r>' +
    '<w:rPr><w:rtl w:val="0"/></w:rPr>' +
    '<w:t>https://hackernoon.com/learning-flask-being-an-ios-developer-3c6ec8c2ba83</w:t>' +
  '</w:r>' +
'</w:p>'

$newPara.Range.InsertXML($newParaXml)
